$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "elapsed_norm" parameter row ---
# Insert a blank row at row 16 (shifts old rows 16-29 down to 17-30, so all
# the other existing parameter rows keep their content/format untouched).
$ws.Rows("16:16").Insert()

# Give the new row 16 the same look (font/alignment/number format) as the
# row directly above it (row 15), which is an existing, fully-populated
# parameter row with the same three-column layout.
$ws.Range("B15:D15").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new parameter: name, description, value.
$ws.Cells.Item(16, 2).Value = "elapsed_norm"
$ws.Cells.Item(16, 3).Value = "elapsed_time normalisation factor"
$ws.Cells.Item(16, 4).Value = 100

# --- Tweak several existing normalisation-factor values (rows above the
# inserted row, so their row numbers are unaffected by the insert) ---
$ws.Cells.Item(8, 4).Value = 0.085000000000000006   # wait_ped_reward_coef
$ws.Cells.Item(10, 4).Value = 45                      # queue_norm
$ws.Cells.Item(11, 4).Value = 100                     # wait_norm
$ws.Cells.Item(12, 4).Value = 7                       # wave_norm
$ws.Cells.Item(13, 4).Value = 20                      # avg_speed_norm
$ws.Cells.Item(15, 4).Value = 100                     # ped_wait_norm

# Match the author's final cursor position.
$ws.Range("D8").Select()
